$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert new customer "Javier Garcia Vaquerizo" as the new row 2
# (pushes every existing data row down by one)
$ws.Rows.Item(2).Insert()
$ws.Range("A2").Value = "021245789C"
$ws.Range("B2").Value = "Javier"
$ws.Range("C2").Value = "Garcia Vaquerizo"
$ws.Range("D2").Value = 658775112
$ws.Range("E2").Value = "C/ Cerveza un Lunes 6"
$ws.Range("F2").Value = "'"
$ws.Range("G2").Value = "ey@yahoo.es"

# Insert new customer "Hector Barrios" as row 19
# (what is currently old row 18 - "26565654C" Alberto - sits at row 19
#  right now, so inserting here pushes it and everything after down by one)
$ws.Rows.Item(19).Insert()
$ws.Range("A19").Value = "254646C"
$ws.Range("B19").Value = "Hector"
$ws.Range("C19").Value = "Barrios  2"
$ws.Range("D19").Value = 65842597
$ws.Range("E19").Value = "C/ Inventada 3"
$ws.Range("F19").Value = '$2a$10$YnzIJs/7gtsOqwz6MJizr.aG1V8bQTjlyhEM3jxGOGL3r5Ko2x8Wm'
$ws.Range("G19").Value = "hector@gmail.com"
